$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Пример")
$ws.Range("B10").Value = "www.stat.gov.kg"
$ws.Activate()
$ws.Range("B10").Select()
